# The test data had two copies of the non-admin user's credentials formula
# stored under column D (duplicating the admin lookup that belongs there),
# and the "Non Admin" rows (4 and 5) were mistakenly still pointing at the
# admin credentials in column C. This corrects the non-admin rows to use
# the proper user_credentials lookup and removes the stray duplicate column
# D formulas on the admin rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transmittals_Close_Cancel")

# Remove the erroneous duplicate "user_credentials" lookup in column D for
# the admin rows (2 and 3).
$ws.Range("D2").ClearContents()
$ws.Range("D3").ClearContents()

# Rows 4 and 5 are the "Non Admin" rows; column C should reference the
# non-admin user (row 4 of user_credentials), not the admin (row 3).
$ws.Range("C4").Formula = "=[1]user_credentials!`$B`$4"
$ws.Range("C5").Formula = "=[1]user_credentials!`$B`$4"
